$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of trade data at row 5
$ws.Range("A5").Value = 9992.61
$ws.Range("B5").Value = 9943.8799999999992
$ws.Range("C5").Value = 307.20999999999998
$ws.Range("D5").Value = 308.70999999999998
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = 0.49
$ws.Range("G5").Value = 42609.503935185188
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = $true
